$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version bump: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date bump
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was previously blank; now populated
$ws.Range("B9").Value = "Alvearie Team"

# The duplicated "Contact" / "No display for ContactDetail" row (row 10)
# becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# The second duplicate "Contact" row (row 11) is removed entirely,
# shifting all subsequent rows up by one.
$ws.Rows.Item(11).Delete()
